$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.957.53"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.535.17"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "306.53"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("E6").Value = "  +8.16%  "
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.548"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").Value = "37.47"
$ws.Range("E10").Value = "  +3.52%  "
$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "7.77"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "0.113"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "2.924.20"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "2.554.63"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "15.27"
$ws.Range("E16").Value = "  +7.36%  "
$ws.Range("D17").Value = "0.872"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "42.955.77"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "13.13"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'6.50"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "71.71"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "254.19"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -3.17%  "
$ws.Range("D26").Value = "27.49"
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "10.51"
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").Value = "2.35"
$ws.Range("E29").Value = "  +10.53%  "
$ws.Range("D30").Value = "38.98"
$ws.Range("E30").Value = "  +4.96%  "
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("D32").Value = "158.74"
$ws.Range("E32").Value = "  +3.14%  "
$ws.Range("D33").Value = "2.11"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "3.32"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0798"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  -3.74%  "
$ws.Range("D37").Value = "18.54"
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "24.25"
$ws.Range("E39").Value = "  +4.58%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.120"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "2.09"
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "3.91"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "2.050.10"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "86.17"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").Value = "8.98"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("D49").Value = "2.785.48"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "103.66"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.193"
$ws.Range("E51").Value = "  +0.77%  "
